# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
#
# This script reproduces, against the already-open workbook, the edits made
# to the "Estado de Cuenta" (account statement) worksheet:
#   1. Insert a new row at row 17 (pushing the trailing signature block
#      down from rows 21/22 to rows 22/23), mirroring Excel's native
#      "insert a copy of the row above" behaviour so the new row inherits
#      row 16's formatting.
#   2. Fill the newly inserted row with a second worker record (same
#      worker, new "Periodo Mora" 2509) underneath the existing one.
#   3. Update "VALOR MORA" (E11) and "Cant. Periodos" (F13) to reflect the
#      newly added period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 17, copying row 16's formatting down so
#    the new data row looks the same as the existing one (borders, fills,
#    number formats, etc.) and the rows below (the signature block) shift
#    down by one (21->22, 22->23).
$ws.Rows.Item(16).Copy()
$ws.Rows.Item(17).Insert()

# 2) New worker record in row 17 (same worker as row 16, new period 2509).
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1010153118"
$ws.Range("D17").Value = "GUSTAVO ANDRES BUSTILLO"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# 3) Update the summary figures: total overdue amount and period count.
$ws.Range("E11").Value = 89206
$ws.Range("F13").Value = 2
